{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// 1. Append \" Govinda govindha\" to the end of the first paragraph\n//    (\"Arunachala shiva \" -> \"Arunachala shiva  Govinda govindha\").\nconst firstParagraph = paragraphs.items[0];\nfirstParagraph.insertText(\" Govinda govindha\", Word.InsertLocation.end);\nawait context.sync();\n\n// 2. Re-write the third paragraph's text in place. The text itself is\n//    unchanged (\"Om namah shivayya.\") but this normalizes the paragraph\n//    down to a single run (clearing the stale spell-check proofErr marks\n//    that were split across the original multi-run text).\nconst thirdParagraph = paragraphs.items[2];\nconst thirdRange = thirdParagraph.getRange();\nthirdRange.insertText(\"Om namah shivayya.\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Append \" Govinda govindha\" to the end of the first paragraph\n#    (\"Arunachala shiva \" -> \"Arunachala shiva  Govinda govindha\").\n$firstParagraph = $d.Paragraphs(1)\n$firstParagraph.Range.InsertAfter(\" Govinda govindha\")\n\n# 2. Re-write the third paragraph's text in place. The text itself is\n#    unchanged (\"Om namah shivayya.\") but this normalizes the paragraph\n#    down to a single run (clearing the stale spell-check proofErr marks\n#    that were split across the original multi-run text).\n$thirdRange = $d.Paragraphs(3).Range\n$thirdRange.Find.ClearFormatting()\n$thirdRange.Find.Execute(\"Om namah shivayya.\", $false, $false, $false, $false, $false, $true, 1, $false, \"Om namah shivayya.\", 2)\n"}
